$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.370.65"
$ws.Range("E2").Value = "  +1.98%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.389.55"
$ws.Range("E3").Value = "  +1.69%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.85"
$ws.Range("E5").Value = "  +0.73%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "179.73"
$ws.Range("E6").Value = "  +1.18%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("E8").Value = "  +0.98%  "

$ws.Range("E9").Value = "  +5.29%  "

$ws.Range("E10").Value = "  +1.41%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "48.47"
$ws.Range("E11").Value = "  +2.57%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000281"
$ws.Range("E12").Value = "  +2.67%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "677.25"
$ws.Range("E13").Value = "  -2.05%  "

$ws.Range("E14").Value = "  +2.24%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.933.57"
$ws.Range("E15").Value = "  +1.66%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "69.430.77"
$ws.Range("E16").Value = "  +2.12%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.396.29"
$ws.Range("E17").Value = "  +1.76%  "

$ws.Range("E18").Value = "  +1.77%  "

$ws.Range("E19").Value = "  +0.94%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.26"
$ws.Range("E20").Value = "  +1.59%  "

$ws.Range("E21").Value = "  +0.50%  "

$ws.Range("E22").Value = "  +0.25%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.13"
$ws.Range("E23").Value = "  +0.39%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "103.39"
$ws.Range("E24").Value = "  +3.79%  "

$ws.Range("E25").Value = "  +0.19%  "

$ws.Range("E26").Value = "  +1.00%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.67"
$ws.Range("E27").Value = "  +0.79%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "34.09"
$ws.Range("E28").Value = "  +2.53%  "

$ws.Range("E29").Value = "  +1.42%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.97"
$ws.Range("E30").Value = "  -2.15%  "

$ws.Range("E31").Value = "  +1.22%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "555.50"
$ws.Range("E32").Value = "  -2.11%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.60"
$ws.Range("E33").Value = "  +6.27%  "

$ws.Range("E34").Value = "  +0.71%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "58.09"
$ws.Range("E35").Value = "  +1.36%  "

$ws.Range("E36").Value = "  -0.02%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.682.71"
$ws.Range("E37").Value = "  -0.76%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.140"
$ws.Range("E38").Value = "  +5.46%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.07"
$ws.Range("E39").Value = "  +1.26%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.28"
$ws.Range("E40").Value = "  +2.55%  "

$ws.Range("E41").Value = "  +0.64%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0₃0698"
$ws.Range("E42").Value = "  +2.97%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.338"
$ws.Range("E43").Value = "  +0.31%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0422"
$ws.Range("E44").Value = "  +3.81%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.31"
$ws.Range("E45").Value = "  -1.10%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.68"
$ws.Range("E46").Value = "  +0.38%  "

$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.42"
$ws.Range("E47").Value = "  +6.42%  "

$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.130"
$ws.Range("E48").Value = "  +0.59%  "

$ws.Range("E49").Value = "  -0.02%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.82"
$ws.Range("E50").Value = "  +1.81%  "

$ws.Range("E51").Value = "  +2.75%  "
